$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates by worksheet and row number, as seen in the diff.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 7242
$ws1.Cells.Item(3, 6).Value = 61
$ws1.Cells.Item(5, 6).Value = 168
$ws1.Cells.Item(6, 6).Value = 1104
$ws1.Cells.Item(7, 6).Value = 182
$ws1.Cells.Item(8, 6).Value = 9
$ws1.Cells.Item(9, 6).Value = 82
$ws1.Cells.Item(10, 6).Value = 19

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 7242
$ws4.Cells.Item(3, 6).Value = 61
$ws4.Cells.Item(5, 6).Value = 168
$ws4.Cells.Item(6, 6).Value = 1104
$ws4.Cells.Item(7, 6).Value = 182
$ws4.Cells.Item(9, 6).Value = 9
$ws4.Cells.Item(10, 6).Value = 82
$ws4.Cells.Item(11, 6).Value = 19

$wb.Save()
